$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Add new rows 61-72 of codebook data ---
# Values must be assigned in this exact order so that the shared-string
# table indices come out in the same sequence as the target workbook.

$ws.Range("A61").Value = "superdomaDiff"
$ws.Range("A62").Value = "domaDiff"
$ws.Range("A63").Value = "trans_disDiff"
$ws.Range("A64").Value = "gay_discDiff"

$ws.Range("B61").Value = "Diffusion variable among the states for superDoma, unlagged"
$ws.Range("B62").Value = "Diffusion variable among the states for doma unlagged"
$ws.Range("B63").Value = "Diffusion variable among the states for transgender antidiscrimination statute (from trans_dis), unlagged"
$ws.Range("B64").Value = "Diffusion variable among the states for gay antidiscrimination statute (from gay_disc), unlagged"

$ws.Range("B66").Value = "lagged version of timeVar by one year"
$ws.Range("B68").Value = "lagged version of timeVarDt by one year"
$ws.Range("B65").Value = "the time varying columns in `"ssph over time with williams.xlsx`" converted to long format and merged with the rest of the data"
$ws.Range("B70").Value = "lagged version of timeVarDtWl by one year"
$ws.Range("B72").Value = "lagged version of timeVarWill by one year"
$ws.Range("B67").Value = "the time varying columns in `"ssphh over time with williams  straight line 1990 2008.xlsx`", converted to long format and merged"
$ws.Range("B69").Value = "the time varying computed williams measures in  `"ssphh over time with williams  straight line 1990 2008.xlsx`", converted to long format and merged"
$ws.Range("B71").Value = "time varying computed williams measures in `"ssphh over time with williams.xlsx`" converted to long format and merged"

$ws.Range("A72").Value = "timeVarWillLag"
$ws.Range("A66").Value = "timeVarLag"
$ws.Range("A67").Value = "timeVarDt"
$ws.Range("A68").Value = "timeVarDtLag"
$ws.Range("A69").Value = "timeVarDtWill"
$ws.Range("A70").Value = "timeVarDtWillLag"
$ws.Range("A71").Value = "timeVarWill"
$ws.Range("A65").Value = "timeVar"

# --- Formatting ---
# B63/B64: plain black font (no hyperlink-like styling)
$ws.Range("B63:B64").Font.Color = 0

# A65:A72: left horizontal alignment
$ws.Range("A65:A72").HorizontalAlignment = -4131

# --- Column sizing ---
$ws.Columns.Item(1).ColumnWidth = 21.33

# --- Page setup ---
$ws.PageSetup.Orientation = 1

# --- View state: scroll so row 54 is at top, select A71 ---
$excel.ActiveWindow.ScrollRow = 54
$ws.Range("A71").Select() | Out-Null
